# Fruta / hortaliza, semanal
# Inserts a new weekly price record as row 90, pushing the existing
# rows 90-107 down to 91-108 (dimension grows from A1:T107 to A1:T108).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 90, shifting
# rows 90:107 down to 91:108.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new data record.
$ws.Range("A90").Value = 5
$ws.Range("B90").Value = "Macroferia Regional de Talca"
$ws.Range("C90").Value = "Maule"
$ws.Range("D90").Value = 44641
$ws.Range("E90").Value = 7
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100103
$ws.Range("H90").Value = "Frutos de hueso (carozo)"
$ws.Range("I90").Value = 100103002
$ws.Range("J90").Value = "Ciruela"
$ws.Range("K90").Value = "Angeleno"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 250
$ws.Range("N90").Value = 8000
$ws.Range("O90").Value = 8000
$ws.Range("P90").Value = 8000
$ws.Range("Q90").Value = "$/bandeja 18 kilos granel"
$ws.Range("R90").Value = "Provincia de Curicó"
$ws.Range("S90").Value = 444
$ws.Range("T90").Value = 18
